$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The second data record (row 3) was re-keyed to become the first (row 2)
# and vice versa: swap the per-record columns (Id, Taxonsorteringsordning,
# TaxonId, Artnamn, Vetenskapligt namn, Auktor, Ost, Nord) between row 2
# and row 3. The remaining columns already hold identical values in both
# rows, so they do not need to change.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "${col}2"
    $addr3 = "${col}3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}

# Column I ("Antal") holds numeric-looking text ("10" / "1") that must stay
# stored as text rather than turn into a real number when swapped.
$i2 = $ws.Range("I2").Value2
$i3 = $ws.Range("I3").Value2
$ws.Range("I2").Value2 = "'" + $i3
$ws.Range("I2").Style = "Normal"
$ws.Range("I3").Value2 = "'" + $i2
$ws.Range("I3").Style = "Normal"
